# Apply the hourly cryptos-list refresh (prices + 1h volume deltas),
# including the Maker/Filecoin row swap at rows 43-44, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.454.67"
$ws.Range("E2").Value = "'  +1.11%  "
$ws.Range("D3").Value = "'3.486.91"
$ws.Range("E3").Value = "'  -0.02%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'594.87"
$ws.Range("E5").Value = "'  +0.32%  "
$ws.Range("D6").Value = "'179.27"
$ws.Range("E6").Value = "'  +4.31%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'3.489.78"
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "'  -0.60%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "'  +5.69%  "
$ws.Range("D11").Value = "'7.09"
$ws.Range("E11").Value = "'  -2.32%  "
$ws.Range("D12").Value = "'0.432"
$ws.Range("E12").Value = "'  +0.06%  "
$ws.Range("D13").Value = "'4.094.48"
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("D14").Value = "'32.24"
$ws.Range("E14").Value = "'  +12.22%  "
$ws.Range("E15").Value = "'  +1.47%  "
$ws.Range("D16").Value = "'67.473.73"
$ws.Range("E16").Value = "'  +1.05%  "
$ws.Range("D17").Value = "'0.0000177"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("D18").Value = "'3.492.88"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("D19").Value = "'6.26"
$ws.Range("E19").Value = "'  -0.25%  "
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "'  +1.96%  "
$ws.Range("D21").Value = "'391.13"
$ws.Range("E21").Value = "'  -0.34%  "
$ws.Range("D22").Value = "'7.88"
$ws.Range("E22").Value = "'  -1.18%  "
$ws.Range("D23").Value = "'72.99"
$ws.Range("E23").Value = "'  +0.46%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "'  -0.17%  "
$ws.Range("D25").Value = "'0.536"
$ws.Range("E25").Value = "'  +0.33%  "
$ws.Range("D26").Value = "'5.71"
$ws.Range("E26").Value = "'  +0.24%  "
$ws.Range("E27").Value = "'  +1.02%  "
$ws.Range("D28").Value = "'10.39"
$ws.Range("E28").Value = "'  +2.14%  "
$ws.Range("D29").Value = "'0.176"
$ws.Range("E29").Value = "'  -2.92%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  +0.25%  "
$ws.Range("E31").Value = "'  -0.34%  "
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("E33").Value = "'  +0.45%  "
$ws.Range("D34").Value = "'23.58"
$ws.Range("E34").Value = "'  -0.14%  "
$ws.Range("D35").Value = "'7.36"
$ws.Range("E35").Value = "'  +0.35%  "
$ws.Range("D37").Value = "'1.60"
$ws.Range("E37").Value = "'  -2.21%  "
$ws.Range("D38").Value = "'164.18"
$ws.Range("E38").Value = "'  +0.58%  "
$ws.Range("D39").Value = "'0.873"
$ws.Range("E39").Value = "'  -0.44%  "
$ws.Range("D40").Value = "'1.87"
$ws.Range("E40").Value = "'  -0.22%  "
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "'  +7.47%  "
$ws.Range("D42").Value = "'6.86"
$ws.Range("E42").Value = "'  -0.81%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.844.99"
$ws.Range("E43").Value = "'  +1.49%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.63"
$ws.Range("E44").Value = "'  -0.68%  "
$ws.Range("D45").Value = "'26.09"
$ws.Range("E45").Value = "'  -0.11%  "
$ws.Range("D46").Value = "'0.0724"
$ws.Range("E46").Value = "'  -2.28%  "
$ws.Range("D47").Value = "'26.50"
$ws.Range("E47").Value = "'  -2.79%  "
$ws.Range("D48").Value = "'42.02"
$ws.Range("E48").Value = "'  -1.20%  "
$ws.Range("D49").Value = "'0.0299"
$ws.Range("E49").Value = "'  -0.93%  "
$ws.Range("D50").Value = "'337.45"
$ws.Range("E50").Value = "'  +0.21%  "
$ws.Range("D51").Value = "'1.05"
$ws.Range("E51").Value = "'  -1.76%  "
